# Update the crypto symbol list (Price / Volume(1h) columns) with freshly
# scraped values, per the GitHub Actions run on Thu Feb 2 11:49:08 UTC 2023.
#
# NOTE: column D (Price) and column E (Volume(1h)) are stored as literal
# text in this sheet (e.g. "330.43", "7.20%") rather than numbers, so every
# assignment below is prefixed with a leading apostrophe to force Excel to
# keep storing them as text instead of auto-converting to numeric/percentage
# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "'329.85"
$ws.Range("E2").Value  = "'7.13%"

$ws.Range("E3").Value  = "'7.65%"

$ws.Range("D4").Value  = "'5.393"
$ws.Range("E4").Value  = "'5.31%"

$ws.Range("D5").Value  = "'0.08138"
$ws.Range("E5").Value  = "'3.97%"

$ws.Range("D6").Value  = "'4.529"

$ws.Range("D7").Value  = "'8.656"
$ws.Range("E7").Value  = "'4.90%"

$ws.Range("D8").Value  = "'1.922"
$ws.Range("E8").Value  = "'2.23%"

$ws.Range("D10").Value = "'0.9439"
$ws.Range("E10").Value = "'2.49%"

$ws.Range("D11").Value = "'0.1364"
$ws.Range("E11").Value = "'26.27%"

$ws.Range("D12").Value = "'0.1979"
$ws.Range("E12").Value = "'4.80%"

$ws.Range("D13").Value = "'0.09336"
$ws.Range("E13").Value = "'5.59%"

$ws.Range("D14").Value = "'0.03554"
$ws.Range("E14").Value = "'7.34%"

$ws.Range("D15").Value = "'0.09590"
$ws.Range("E15").Value = "'-0.02%"

$ws.Range("D16").Value = "'0.001322"
$ws.Range("E16").Value = "'-3.99%"

$ws.Range("E17").Value = "'10.01%"

$ws.Range("D18").Value = "'3.362"
$ws.Range("E18").Value = "'-1.13%"

$ws.Range("E19").Value = "'2.89%"

$ws.Range("D20").Value = "'7.209"
$ws.Range("E20").Value = "'14.35%"

$ws.Range("E21").Value = "'2.10%"

$ws.Range("D22").Value = "'0.2557"
$ws.Range("E22").Value = "'5.88%"

$ws.Range("D23").Value = "'0.04432"
$ws.Range("E23").Value = "'1.32%"

$ws.Range("D24").Value = "'0.001223"
$ws.Range("E24").Value = "'2.40%"

$ws.Range("D25").Value = "'0.004296"
$ws.Range("E25").Value = "'0.67%"

$ws.Range("E26").Value = "'-14.28%"

$ws.Range("E27").Value = "'-0.05%"

$ws.Range("D39").Value = "'0.02497"
$ws.Range("E39").Value = "'15.00%"

$ws.Range("D40").Value = "'0.05234"
$ws.Range("E40").Value = "'3.76%"

$ws.Range("D41").Value = "'0.007584"
$ws.Range("E41").Value = "'0.45%"

$ws.Range("D42").Value = "'0.1428"
$ws.Range("E42").Value = "'5.57%"

$ws.Range("D43").Value = "'0.009131"
$ws.Range("E43").Value = "'5.47%"

$ws.Range("D44").Value = "'0.002170"
$ws.Range("E44").Value = "'4.92%"

$ws.Range("E45").Value = "'38.00%"

$ws.Range("D46").Value = "'0.00006590"
$ws.Range("E46").Value = "'1.16%"

$ws.Range("E47").Value = "'0.00%"

$ws.Range("E48").Value = "'139.48%"

$ws.Range("E49").Value = "'1.55%"

$ws.Range("E50").Value = "'0.00%"

$ws.Range("E51").Value = "'0.00%"
